$d = $word.ActiveDocument

# --- Step 1: remove the old "_GoBack" bookmark that currently sits mid-run
#     in the "commit" paragraph (it will be re-created after the new notes
#     are appended, mirroring Word's own "last edit" bookmark behaviour). ---
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- Step 2: append the two new note paragraphs right after the
#     "commit" paragraph, and restore the _GoBack bookmark on the
#     (still) final, empty paragraph. We build these as raw OOXML so the
#     new paragraphs come out "clean" (no inherited <w:pPr>), and splice
#     them in with a single InsertXML call: Word folds the last fragment
#     of inserted XML into the paragraph that follows the insertion
#     point, so by ending our payload with a bookmark-only paragraph we
#     land the bookmark back on the original trailing paragraph without
#     disturbing its existing <w:pPr>. ---
$commitPara = $d.Paragraphs.Item(5)
$insertionPoint = $d.Range($commitPara.Range.End, $commitPara.Range.End)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$statusPara = "<w:p $w>" +
    "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>查看状态：git</w:t></w:r>" +
    "<w:r><w:t xml:space=""preserve""> status</w:t></w:r>" +
    "</w:p>"

$diffPara = "<w:p $w>" +
    "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>查看修改：git</w:t></w:r>" +
    "<w:r><w:t xml:space=""preserve""> </w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>diff</w:t></w:r>" +
    "</w:p>"

$trailingPara = "<w:p $w>" +
    "<w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/>" +
    "</w:p>"

[void]$insertionPoint.InsertXML($statusPara + $diffPara + $trailingPara)
